$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Label the "Tooling" week-2 cell as Linux, and mark it with the
# "Should be done" (red) color used elsewhere in the legend.
$ws.Range("C14").Value = "Linux"
$ws.Range("C14").Interior.Color = $ws.Range("B18").Interior.Color

# Fix capitalization of the legend entries.
$ws.Range("A21").Value = "Worked on that week"
$ws.Range("A20").Value = "Planned to be done that week"

# Update the saved selection on the sheet.
$ws.Range("C4").Select()
